{"js": "// \"Refactoring and adding test cases.\"\n//\n//  1. Remove the `_GoBack` bookmark that currently sits in the Heading1\n//     paragraph (\"Test Specs\").\n//  2. Rewrite the body paragraph's text: expand on the load/unload test\n//     case and drop the trailing sentence about adding/removing plugins.\n//     (This also removes the spell-check `proofErr` markers that used to\n//     bracket \"JUnit\" - purely editorial artifacts that aren't recreated\n//     when the run is rewritten.)\n//  3. Re-insert the `_GoBack` bookmark, now positioned just before the\n//     final \".\" of the rewritten paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst secondParagraph = paragraphs.items[1];\n\n// 1. Drop the existing `_GoBack` bookmark (lives in the heading paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Replace the whole paragraph text with the new wording.\nconst newText =\n  \"We will use unit testing with the JUnit framework to test specific \" +\n  \"methods of importance to ensure they behave as expected.  The methods \" +\n  \"we will test are load and unload bundle in the plugin manager, as \" +\n  \"these are the most important methods in the project.\";\nsecondParagraph.insertText(newText, \"Replace\");\nawait context.sync();\n\n// 3. Re-insert the `_GoBack` bookmark right before the trailing period.\nconst searchResults = secondParagraph.search(\"in the project\", {\n  matchCase: true,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst insertionPoint = searchResults.items[0].getRange(\"End\");\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"Refactoring and adding test cases.\"\n#\n#  1. Remove the `_GoBack` bookmark that currently sits in the Heading1\n#     paragraph (\"Test Specs\").\n#  2. Rewrite the body paragraph's text: expand on the load/unload test\n#     case and drop the trailing sentence about adding/removing plugins.\n#     (This also removes the spell-check proofing marks that used to\n#     bracket \"JUnit\" - purely editorial artifacts that aren't recreated\n#     when the paragraph text is rewritten.)\n#  3. Re-insert the `_GoBack` bookmark, now positioned just before the\n#     final \".\" of the rewritten paragraph.\n\n$d = $word.ActiveDocument\n\n# 1. Drop the existing `_GoBack` bookmark (lives in the heading paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Replace the whole second paragraph's text (but not its paragraph\n#    mark) with the new wording.\n$paraRange = $d.Paragraphs.Item(2).Range\n$bodyRange = $d.Range($paraRange.Start, $paraRange.End - 1)\n$bodyRange.Text = \"We will use unit testing with the JUnit framework to test specific methods of importance to ensure they behave as expected.  The methods we will test are load and unload bundle in the plugin manager, as these are the most important methods in the project.\"\n\n# 3. Re-insert the `_GoBack` bookmark right before the trailing period.\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Execute(\"in the project\") | Out-Null\n$insertionPoint = $d.Range($searchRange.End, $searchRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n"}
